$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "~, 12, 13, 14, 16, 17"
$ws.Range("B5").Value = "~, 12, 13, 14, 16"
$ws.Range("B6").Value = "~, 16, 17"
$ws.Range("B7").Value = "~, 12, 13, 14, 15"
$ws.Range("B8").Value = "~, 16"
$ws.Range("B9").Value = "~, 12, 13, 14, 15"
$ws.Range("B11").Value = "~, 17"
$ws.Range("B15").Value = "~, 14, 16"
$ws.Range("B16").Value = "~, 15, 16, 17"
$ws.Range("B18").Value = "~, 15"
$ws.Range("B22").Value = "~, 14, 17"
$ws.Range("B24").Value = "~, 12, 13, 15"
$ws.Range("B25").Value = "~, 12, 13, 16"
$ws.Range("B26").Value = "~, 14"
$ws.Range("B29").Value = "~, 12, 13, 14, 15, 16, 17"
$ws.Range("B31").Value = "~, 12, 13, 14, 15, 16"
$ws.Range("B32").Value = "~, 16"
$ws.Range("B33").Value = "~, 13, 14, 15, 17"
$ws.Range("B35").Value = "~, 13, 14, 15, 16"
$ws.Range("B36").Value = "~, 15"
$ws.Range("B37").Value = "~, 12, 13, 14"
$ws.Range("B39").Value = "~, 12, 13, 14"
$ws.Range("B41").Value = "~, 12, 13, 14"
$ws.Range("B42").Value = "~, 12, 13, 16"
$ws.Range("B43").Value = "~, 14"
$ws.Range("B45").Value = "~, 13, 14, 16, 17"
$ws.Range("B47").Value = "~, 12, 13, 14, 15"
$ws.Range("B48").Value = "~, 16"
$ws.Range("B49").Value = "~, 12, 13, 14, 15, 17"
$ws.Range("B50").Value = "~, 12, 14"
$ws.Range("B51").Value = "~, 13, 15, 16"
$ws.Range("B53").Value = "~, 12, 15, 16"
$ws.Range("B66").Value = "~, 12, 13, 14, 15, 16"
$ws.Range("D66").Value = "TARDE"
$ws.Range("E66").Value = "ROTA NÃO CONDIZ"
$ws.Range("B67").Value = "~, 12, 13, 14, 15, 16"
$ws.Range("E67").Value = "ROTA NÃO CONDIZ"
$ws.Range("A68").Value = "NQR-5926"
$ws.Range("D68").Value = "INTEGRAL"
$ws.Range("E68").Value = "NÃO APRESENTA ROTA"
$ws.Range("A69").Value = "NQR-5926"
$ws.Range("D69").Value = "INTEGRAL"
$ws.Range("E69").Value = "NÃO FEZ A ROTA"
$ws.Range("A70").Value = "KJK-9345"
$ws.Range("B70").Value = "~, 13"
$ws.Range("F70").Value = "UNIÃO"
$ws.Range("A71").Value = "KJK-9345"
$ws.Range("B71").Value = "~, 12"
$ws.Range("F71").Value = "UNIÃO"
$ws.Range("A72").Value = "KXR-5549"
$ws.Range("B72").Value = "~, 12, 13, 14, 15, 16"
$ws.Range("A73").Value = "KXR-5549"
$ws.Range("A74").Value = "NGJ-0903"
$ws.Range("B74").Value = "~, 12"
$ws.Range("A75").Value = "NGJ-0903"
$ws.Range("B75").Value = "~, 16"
$ws.Range("A76").Value = "NGJ-0903"
$ws.Range("B76").Value = "~, 13, 14"
$ws.Range("C76").Value = "JUNHO"
$ws.Range("D76").Value = "TARDE"
$ws.Range("E76").Value = "ROTA NÃO CONDIZ"
$ws.Range("F76").Value = "PALMERAIS"
$ws.Range("A77").Value = "NGC-8853"
$ws.Range("B77").Value = "~"
$ws.Range("C77").Value = "JUNHO"
$ws.Range("D77").Value = "TARDE"
$ws.Range("E77").Value = "NÃO FEZ A ROTA"
$ws.Range("F77").Value = "PALMERAIS"
$ws.Range("A78").Value = "NGC-8853"
$ws.Range("B78").Value = "~, 12, 13, 14"
$ws.Range("C78").Value = "JUNHO"
$ws.Range("D78").Value = "TARDE"
$ws.Range("E78").Value = "NÃO APRESENTA ROTA"
$ws.Range("F78").Value = "PALMERAIS"

Write-Output "done"